$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.959.20"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").Value = "'1.867.51"
$ws.Range("E3").Value = "  -2.75%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'318.85"
$ws.Range("E5").Value = "  -3.29%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D7").Value = "'0.5083"
$ws.Range("E7").Value = "  -3.33%  "
$ws.Range("D8").Value = "'0.3935"
$ws.Range("E8").Value = "  -2.82%  "
$ws.Range("D9").Value = "'0.08203"
$ws.Range("E9").Value = "  -3.24%  "
$ws.Range("D10").Value = "'42.18"
$ws.Range("E10").Value = "  -1.66%  "
$ws.Range("D11").Value = "'1.092"
$ws.Range("E11").Value = "  -3.36%  "
$ws.Range("D12").Value = "'22.88"
$ws.Range("E12").Value = "  +2.88%  "
$ws.Range("D13").Value = "'1.873.95"
$ws.Range("E13").Value = "  -2.42%  "
$ws.Range("D14").Value = "'6.273"
$ws.Range("E14").Value = "  -1.68%  "
$ws.Range("D15").Value = "'7.168"
$ws.Range("E15").Value = "  -2.98%  "
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "'91.98"
$ws.Range("E17").Value = "  -4.63%  "
$ws.Range("D18").Value = "'0.00001085"
$ws.Range("E18").Value = "  -2.69%  "
$ws.Range("D19").Value = "'0.06381"
$ws.Range("E19").Value = "  -5.00%  "
$ws.Range("D20").Value = "'17.87"
$ws.Range("E20").Value = "  -2.12%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "'29.940.46"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").Value = "'5.818"
$ws.Range("E23").Value = "  -4.06%  "
$ws.Range("D24").Value = "'11.11"
$ws.Range("E24").Value = "  -1.39%  "
$ws.Range("D25").Value = "'2.172"
$ws.Range("E25").Value = "  -2.39%  "
$ws.Range("D26").Value = "'2.081.38"
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("D27").Value = "'161.11"
$ws.Range("E27").Value = "  +0.51%  "
$ws.Range("D28").Value = "'20.93"
$ws.Range("D29").Value = "'2.221"
$ws.Range("E29").Value = "  -9.65%  "
$ws.Range("D30").Value = "'127.45"
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("D31").Value = "'1.062"
$ws.Range("E31").Value = "  -2.05%  "
$ws.Range("E32").Value = "  -2.65%  "
$ws.Range("D33").Value = "'5.907"
$ws.Range("E33").Value = "  -3.39%  "
$ws.Range("D34").Value = "'3.725"
$ws.Range("E34").Value = "  +1.80%  "
$ws.Range("E35").Value = "  -3.90%  "
$ws.Range("D36").Value = "'5.217"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("D37").Value = "'0.06335"
$ws.Range("E37").Value = "  -4.00%  "
$ws.Range("D38").Value = "'0.2138"
$ws.Range("E38").Value = "  -4.36%  "
$ws.Range("D39").Value = "'1.170"
$ws.Range("E39").Value = "  -5.36%  "
$ws.Range("D40").Value = "'8.485"
$ws.Range("E40").Value = "  -6.03%  "
$ws.Range("D41").Value = "'0.6292"
$ws.Range("E41").Value = "  -4.17%  "
$ws.Range("D42").Value = "'1.205"
$ws.Range("E42").Value = "  -3.45%  "
$ws.Range("D43").Value = "'11.25"
$ws.Range("E43").Value = "  -4.37%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").Value = "'0.5887"
$ws.Range("E45").Value = "  -5.08%  "
$ws.Range("D46").Value = "'12.90"
$ws.Range("E46").Value = "  -2.82%  "
$ws.Range("D47").Value = "'3.633"
$ws.Range("E47").Value = "  -3.72%  "
$ws.Range("D48").Value = "'1.999"
$ws.Range("E48").Value = "  -3.79%  "
$ws.Range("D49").Value = "'122.34"
$ws.Range("E49").Value = "  -2.90%  "
$ws.Range("D50").Value = "'1.201"
$ws.Range("E50").Value = "  -3.37%  "
$ws.Range("D51").Value = "'1.121"
$ws.Range("E51").Value = "  -2.84%  "

Write-Host "Updated cryptos list"
